# Adds two new worksheets ("create_single_user_data" and "update_user_details")
# to the PetStoreTestData workbook, describing the "update user details" test
# flow (update firstName / lastName / email / password / phone), plus a
# "create single user" data sheet. Mirrors the target commit:
#   "Flow for the update user details ... Also added the flow for get user
#    details API"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper-ish constants
# ---------------------------------------------------------------------------
$xlCenter = -4108

$headers = @("TestCaseNo","Description","id","username","firstName","lastName","email","password","phone","userStatus","ifToRunTC")

# ===========================================================================
# Sheet 1 (create_users): just drop the explicit "active cell" / tab
# selection left over from the previous save, matching the new view state.
# ===========================================================================
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:K2").Select()

# ===========================================================================
# Sheet 2: create_single_user_data
# ===========================================================================
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "create_single_user_data"

# --- header row -------------------------------------------------------------
for ($col = 1; $col -le 11; $col++) {
    $cell = $ws2.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = $xlCenter
    if ($col -ne 2) {
        $cell.VerticalAlignment = $xlCenter
    }
}

# --- data row ----------------------------------------------------------------
$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(2, 2).Value = "Common_user_details"
$ws2.Cells.Item(2, 3).Value = 111
$ws2.Cells.Item(2, 4).Value = "testUser111"
$ws2.Cells.Item(2, 5).Value = "testFirst111"
$ws2.Cells.Item(2, 6).Value = "testLast111"
$ws2.Cells.Item(2, 7).Value = "testuser111@gmail.com"
$ws2.Cells.Item(2, 8).Value = "test123"
$ws2.Cells.Item(2, 9).Value = 1234567890
$ws2.Cells.Item(2, 10).Value = 0
$ws2.Cells.Item(2, 11).Value = "Y"

# bold the "Common" prefix of the description cell
$ws2.Cells.Item(2, 2).Characters(1, 6).Font.Bold = $true

# border the whole used range, center-align it
$usedRange2 = $ws2.Range("A1:K2")
$usedRange2.Borders.LineStyle = 1
$rng = $ws2.Range("A2:K2")
$rng.HorizontalAlignment = $xlCenter

# hyperlink the e-mail cell
$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 7), "mailto:testuser111@gmail.com")

# column widths (approximate, matches bestFit widths from Excel)
$ws2.Columns.Item(1).ColumnWidth = 10.83
$ws2.Columns.Item(2).ColumnWidth = 19.66
$ws2.Columns.Item(3).ColumnWidth = 12
$ws2.Range("D1:F1").ColumnWidth = 16.83
$ws2.Columns.Item(7).ColumnWidth = 24.33
$ws2.Range("H1:K1").ColumnWidth = 16.83

$ws2.Cells.Item(2, 4).Select()

# ===========================================================================
# Sheet 3: update_user_details
# ===========================================================================
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "update_user_details"

for ($col = 1; $col -le 11; $col++) {
    $cell = $ws3.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = $xlCenter
    if ($col -ne 2) {
        $cell.VerticalAlignment = $xlCenter
    }
}

$rows = @(
    @{ n=1; desc="update_firstName"; boldStart=7;  boldLen=10; email="testuser111@gmail.com";         first="testFirst111_updated"; last="testLast111";          pwd="test123";         phone=1234567890 },
    @{ n=2; desc="update_lastName";  boldStart=7;  boldLen=9;  email="testuser111@gmail.com";         first="testFirst111_updated"; last="testLast111_updated";  pwd="test123";         phone=1234567890 },
    @{ n=3; desc="update_email";     boldStart=7;  boldLen=6;  email="testuser111_updated@gmail.com"; first="testFirst111_updated"; last="testLast111_updated";  pwd="test123";         phone=1234567890 },
    @{ n=4; desc="update_password";  boldStart=7;  boldLen=9;  email="testuser111_updated@gmail.com"; first="testFirst111_updated"; last="testLast111_updated";  pwd="test123_updated"; phone=1234567890 },
    @{ n=5; desc="update_phone";     boldStart=7;  boldLen=6;  email="testuser111_updated@gmail.com"; first="testFirst111_updated"; last="testLast111_updated";  pwd="test123_updated"; phone=1234567899 }
)

$r = 2
foreach ($row in $rows) {
    $ws3.Cells.Item($r, 1).Value = $row.n
    $ws3.Cells.Item($r, 2).Value = $row.desc
    $ws3.Cells.Item($r, 2).Characters($row.boldStart, $row.boldLen).Font.Bold = $true
    $ws3.Cells.Item($r, 3).Value = 111
    $ws3.Cells.Item($r, 4).Value = "testUser111"
    $ws3.Cells.Item($r, 5).Value = $row.first
    $ws3.Cells.Item($r, 6).Value = $row.last
    $ws3.Cells.Item($r, 7).Value = $row.email
    $ws3.Cells.Item($r, 8).Value = $row.pwd
    $ws3.Cells.Item($r, 9).Value = $row.phone
    $ws3.Cells.Item($r, 10).Value = 0
    $ws3.Cells.Item($r, 11).Value = "Y"
    $r++
}

# border + center whole used range
$usedRange3 = $ws3.Range("A1:K6")
$usedRange3.Borders.LineStyle = 1
$ws3.Range("A2:K6").HorizontalAlignment = $xlCenter

# highlight (yellow) the cells that changed value versus the previous row,
# i.e. the "what this scenario updates" cells
$yellow = 65535
$ws3.Range("E2").Interior.Color = $yellow         # updated firstName (row 1)
$ws3.Range("F3").Interior.Color = $yellow         # updated lastName (row 2)
$ws3.Range("G4").Interior.Color = $yellow         # updated email (row 3)
$ws3.Range("H5").Interior.Color = $yellow         # updated password (row 4)
$ws3.Range("I6").Interior.Color = $yellow         # updated phone (row 5)

# hyperlink the e-mail cells
$ws3.Hyperlinks.Add($ws3.Cells.Item(2, 7), "mailto:testuser111@gmail.com")
$ws3.Hyperlinks.Add($ws3.Cells.Item(3, 7), "mailto:testuser111@gmail.com")
$ws3.Hyperlinks.Add($ws3.Cells.Item(4, 7), "mailto:testuser111_updated@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("G5:G6"), "mailto:testuser111_updated@gmail.com", "", "", "testuser111_updated@gmail.com")

# column widths
$ws3.Columns.Item(2).ColumnWidth = 19.66
$ws3.Columns.Item(3).ColumnWidth = 9.66
$ws3.Columns.Item(4).ColumnWidth = 16.5
$ws3.Columns.Item(5).ColumnWidth = 18.83
$ws3.Columns.Item(6).ColumnWidth = 20.5
$ws3.Columns.Item(7).ColumnWidth = 29.5
$ws3.Columns.Item(8).ColumnWidth = 15
$ws3.Columns.Item(9).ColumnWidth = 14.83

$ws3.Cells.Item(12, 8).Select()

# the last-edited / active sheet is "update_user_details"
$ws3.Activate()
